$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "book_value" / "market_value" row formatting (label cell in col A,
# and the "0" filler cells in cols B:D) down onto the two new rows, exactly
# like a user dragging/copy-pasting the existing rows before retyping them.
$ws.Range("A8:D9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 10: book_value_return ---
$ws.Range("A10").Value = "book_value_return"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0.010273972602739725
$ws.Range("F10").Value = 0.03962575674188222
$ws.Range("G10").Value = -0.03139717425431711
$ws.Range("H10").Value = 0.02592352559948153
$ws.Range("I10").Value = 0.030749100425253518
$ws.Range("J10").Value = 0.0386310604096448
$ws.Range("K10").Value = 0.0487651077246453
$ws.Range("L10").Value = 0.062218959346609155
$ws.Range("M10").Value = 0.05001190759704691
$ws.Range("N10").Value = 0.03855109961190168
$ws.Range("O10").Value = 0.04431609346285114
$ws.Range("P10").Value = 0.035305711470071555
$ws.Range("Q10").Value = 0.13207636553128996
$ws.Range("R10").Value = 0.05397859469520707
$ws.Range("S10").Value = 0.07285301501586851

# --- Row 11: market_value_return ---
$ws.Range("A11").Value = "market_value_return"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0.056
$ws.Range("F11").Value = -0.224
$ws.Range("G11").Value = 0.113
$ws.Range("H11").Value = 0.123
$ws.Range("I11").Value = 0.04
$ws.Range("J11").Value = 0.084
$ws.Range("K11").Value = 0.12
$ws.Range("L11").Value = 0.117
$ws.Range("M11").Value = 0.081
$ws.Range("N11").Value = 0.073
$ws.Range("O11").Value = 0.085
$ws.Range("P11").Value = 0.044
$ws.Range("Q11").Value = 0.095
$ws.Range("R11").Value = 0.07
$ws.Range("S11").Value = 0.114

# The two new "return" rows get a dedicated 4-decimal number format
# (book value is derived via linear-regression coefficients; market value
# return uses a 2-year lag), distinct from the integer/thousands formats
# used by the amount rows above.
$ws.Range("E10:S11").NumberFormat = "0.0000"

# Match the saved selection/window state left behind after entering the data.
$ws.Range("D15").Select()
